# Export with no is_pref and no lev distance
#
# Re-generates the "id" (column B) and "speaker_variant" (column C) values for
# every data row of the playlist/speaker-variant sheet, and clears the
# "is_prefered" flag (column D) that used to mark a single preferred variant
# per speaker. The new export no longer picks a preferred variant (is_pref)
# nor collapses near-duplicate spellings using Levenshtein distance, so every
# distinct spelling variant now shows up as its own row with D left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = '#p.-paulina'
$ws.Cells.Item(2, 3).Value = 'P. Paulina'
$ws.Cells.Item(2, 4).Value = $null
$ws.Cells.Item(3, 2).Value = '#faon'
$ws.Cells.Item(3, 3).Value = 'Faon'
$ws.Cells.Item(3, 4).Value = $null
$ws.Cells.Item(4, 2).Value = '#epaphrod'
$ws.Cells.Item(4, 3).Value = 'Epaphrod'
$ws.Cells.Item(4, 4).Value = $null
$ws.Cells.Item(5, 2).Value = '#tigellinus'
$ws.Cells.Item(5, 3).Value = 'Tigellinus'
$ws.Cells.Item(5, 4).Value = $null
$ws.Cells.Item(6, 2).Value = '#tigellinu'
$ws.Cells.Item(6, 3).Value = 'Tigellinu'
$ws.Cells.Item(6, 4).Value = $null
$ws.Cells.Item(7, 2).Value = '#sporus'
$ws.Cells.Item(7, 3).Value = 'Sporus'
$ws.Cells.Item(7, 4).Value = $null
$ws.Cells.Item(8, 2).Value = '#burrus'
$ws.Cells.Item(8, 3).Value = 'Burrus'
$ws.Cells.Item(8, 4).Value = $null
$ws.Cells.Item(9, 2).Value = '#praetorian'
$ws.Cells.Item(9, 3).Value = 'Praetorian'
$ws.Cells.Item(9, 4).Value = $null
$ws.Cells.Item(10, 2).Value = '#roma,-tyber,-ende-choor-singhen,-oft-spreken-te-samen'
$ws.Cells.Item(10, 3).Value = 'Roma, Tyber, ende Choor singhen, oft spreken te samen'
$ws.Cells.Item(10, 4).Value = $null
$ws.Cells.Item(11, 2).Value = '#eccl.-voest'
$ws.Cells.Item(11, 3).Value = 'Eccl. Voest'
$ws.Cells.Item(11, 4).Value = $null
$ws.Cells.Item(12, 2).Value = '#k.-va'
$ws.Cells.Item(12, 3).Value = 'k. va'
$ws.Cells.Item(12, 4).Value = $null
$ws.Cells.Item(13, 2).Value = '#megara'
$ws.Cells.Item(13, 3).Value = 'Megara'
$ws.Cells.Item(13, 4).Value = $null
$ws.Cells.Item(14, 2).Value = '#agrippina'
$ws.Cells.Item(14, 3).Value = 'Agrippina'
$ws.Cells.Item(14, 4).Value = $null
$ws.Cells.Item(15, 2).Value = '#herault'
$ws.Cells.Item(15, 3).Value = 'Herault'
$ws.Cells.Item(15, 4).Value = $null
$ws.Cells.Item(16, 2).Value = '#pretorius'
$ws.Cells.Item(16, 3).Value = 'Pretorius'
$ws.Cells.Item(16, 4).Value = $null
$ws.Cells.Item(17, 2).Value = '#octauia'
$ws.Cells.Item(17, 3).Value = 'Octauia'
$ws.Cells.Item(17, 4).Value = $null
$ws.Cells.Item(18, 2).Value = '#fenius'
$ws.Cells.Item(18, 3).Value = 'Fenius'
$ws.Cells.Item(18, 4).Value = $null
$ws.Cells.Item(19, 2).Value = '#epaphrodi'
$ws.Cells.Item(19, 3).Value = 'Epaphrodi'
$ws.Cells.Item(19, 4).Value = $null
$ws.Cells.Item(20, 2).Value = '#c.-senici'
$ws.Cells.Item(20, 3).Value = 'C. Senici'
$ws.Cells.Item(20, 4).Value = $null
$ws.Cells.Item(21, 2).Value = '#siluanus'
$ws.Cells.Item(21, 3).Value = 'Siluanus'
$ws.Cells.Item(21, 4).Value = $null
$ws.Cells.Item(22, 2).Value = '#icelus'
$ws.Cells.Item(22, 3).Value = 'Icelus'
$ws.Cells.Item(22, 4).Value = $null
$ws.Cells.Item(23, 2).Value = '#siluanus-alleen'
$ws.Cells.Item(23, 3).Value = 'Siluanus alleen'
$ws.Cells.Item(23, 4).Value = $null
$ws.Cells.Item(24, 2).Value = '#laquay'
$ws.Cells.Item(24, 3).Value = 'Laquay'
$ws.Cells.Item(24, 4).Value = $null
$ws.Cells.Item(25, 2).Value = '#te-samen-als-te-vooren'
$ws.Cells.Item(25, 3).Value = 'Te samen als te vooren'
$ws.Cells.Item(25, 4).Value = $null
$ws.Cells.Item(26, 2).Value = '#acté'
$ws.Cells.Item(26, 3).Value = 'Acté'
$ws.Cells.Item(26, 4).Value = $null
$ws.Cells.Item(27, 2).Value = '#ecloge-voester'
$ws.Cells.Item(27, 3).Value = 'Ecloge Voester'
$ws.Cells.Item(27, 4).Value = $null
$ws.Cells.Item(28, 2).Value = '#statius'
$ws.Cells.Item(28, 3).Value = 'Statius'
$ws.Cells.Item(28, 4).Value = $null
$ws.Cells.Item(29, 2).Value = '#tyber'
$ws.Cells.Item(29, 3).Value = 'Tyber'
$ws.Cells.Item(29, 4).Value = $null
$ws.Cells.Item(30, 2).Value = '#nymphid'
$ws.Cells.Item(30, 3).Value = 'Nymphid'
$ws.Cells.Item(30, 4).Value = $null
$ws.Cells.Item(31, 2).Value = '#anicetu'
$ws.Cells.Item(31, 3).Value = 'Anicetu'
$ws.Cells.Item(31, 4).Value = $null
$ws.Cells.Item(32, 2).Value = '#nymphidi'
$ws.Cells.Item(32, 3).Value = 'Nymphidi'
$ws.Cells.Item(32, 4).Value = $null
$ws.Cells.Item(33, 2).Value = '#mellichus'
$ws.Cells.Item(33, 3).Value = 'Mellichus'
$ws.Cells.Item(33, 4).Value = $null
$ws.Cells.Item(34, 2).Value = '#c.-seneci'
$ws.Cells.Item(34, 3).Value = 'C. Seneci'
$ws.Cells.Item(34, 4).Value = $null
$ws.Cells.Item(35, 2).Value = '#petinu'
$ws.Cells.Item(35, 3).Value = 'Petinu'
$ws.Cells.Item(35, 4).Value = $null
$ws.Cells.Item(36, 2).Value = '#nero'
$ws.Cells.Item(36, 3).Value = 'Nero'
$ws.Cells.Item(36, 4).Value = $null
$ws.Cells.Item(37, 2).Value = '#anicetus'
$ws.Cells.Item(37, 3).Value = 'Anicetus'
$ws.Cells.Item(37, 4).Value = $null
$ws.Cells.Item(38, 2).Value = '#ecl'
$ws.Cells.Item(38, 3).Value = 'Ecl'
$ws.Cells.Item(38, 4).Value = $null
$ws.Cells.Item(39, 2).Value = '#epilogus'
$ws.Cells.Item(39, 3).Value = 'Epilogus'
$ws.Cells.Item(39, 4).Value = $null
$ws.Cells.Item(40, 2).Value = '#seneca'
$ws.Cells.Item(40, 3).Value = 'Seneca'
$ws.Cells.Item(40, 4).Value = $null
$ws.Cells.Item(41, 2).Value = '#s.-poppea'
$ws.Cells.Item(41, 3).Value = 'S. Poppea'
$ws.Cells.Item(42, 2).Value = '#soldae'
$ws.Cells.Item(42, 3).Value = 'Soldae'
$ws.Cells.Item(43, 2).Value = '#soldaet'
$ws.Cells.Item(43, 3).Value = 'Soldaet'
$ws.Cells.Item(44, 2).Value = '#choor'
$ws.Cells.Item(44, 3).Value = 'CHOOR'
$ws.Cells.Item(45, 2).Value = '#a.-natalis'
$ws.Cells.Item(45, 3).Value = 'A. Natalis'
$ws.Cells.Item(46, 2).Value = '#roma'
$ws.Cells.Item(46, 3).Value = 'Roma'
$ws.Cells.Item(47, 2).Value = '#ecl.-voest'
$ws.Cells.Item(47, 3).Value = 'Ecl. Voest'
